$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(297, "Login with valid username and password", "PASSED", "chrome"),
    @(298, "Create a Citizenship", "FAILED", "chrome"),
    @(299, "Login with valid username and password", "PASSED", "chrome"),
    @(300, "Create a country", "PASSED", "chrome"),
    @(301, "Create a country", "PASSED", "chrome"),
    @(302, "Create a country", "FAILED", "chrome"),
    @(303, "Create a country", "PASSED", "chrome"),
    @(304, "Create a country", "PASSED", "chrome"),
    @(305, "Create a country", "PASSED", "chrome"),
    @(306, "Create a country", "PASSED", "chrome"),
    @(307, "Create a country", "PASSED", "chrome"),
    @(308, "Create a country", "PASSED", "chrome"),
    @(309, "Create a country", "PASSED", "chrome"),
    @(310, "Create a country 2", "PASSED", "chrome"),
    @(311, "Create a country", "PASSED", "chrome"),
    @(312, "Create a country 2", "PASSED", "chrome"),
    @(313, "Create a citizenship", "PASSED", "chrome"),
    @(314, "Create a citizenship", "FAILED", "chrome"),
    @(315, "Create a citizenship", "PASSED", "chrome"),
    @(316, "Create a citizenship", "FAILED", "chrome"),
    @(317, "Create a citizenship", "PASSED", "chrome"),
    @(318, "Create a citizenship", "PASSED", "chrome"),
    @(319, "Create a citizenship", "PASSED", "chrome"),
    @(320, "Create a citizenship", "PASSED", "chrome"),
    @(321, "Create a citizenship", "PASSED", "chrome"),
    @(322, "Create a citizenship", "PASSED", "chrome"),
    @(323, "Create a Citizenship", "PASSED", "chrome"),
    @(324, "Create a Citizenship", "PASSED", "chrome"),
    @(325, "Create a Citizenship", "PASSED", "chrome"),
    @(326, "Create a Citizenship", "PASSED", "chrome"),
    @(327, "Create a Citizenship", "PASSED", "chrome"),
    @(328, "Country with Parameter", "FAILED", "chrome"),
    @(329, "Country with Parameter", "FAILED", "chrome"),
    @(330, "Country with Parameter", "FAILED", "chrome"),
    @(331, "Country with Parameter", "FAILED", "chrome"),
    @(332, "Country with Parameter", "FAILED", "chrome"),
    @(333, "Country with Parameter", "FAILED", "chrome"),
    @(334, "Create Country", "FAILED", "chrome"),
    @(335, "Create Country", "FAILED", "chrome"),
    @(336, "Create Country", "PASSED", "chrome"),
    @(337, "Create Nationality", "PASSED", "chrome"),
    @(338, "Fee Functionality", "FAILED", "chrome"),
    @(339, "Create Country", "PASSED", "chrome"),
    @(340, "Create Nationality", "PASSED", "chrome"),
    @(341, "Fee Functionality", "FAILED", "chrome"),
    @(342, "Create Country", "PASSED", "chrome"),
    @(343, "Create Nationality", "PASSED", "chrome"),
    @(344, "Fee Functionality", "FAILED", "chrome"),
    @(345, "Create Country", "PASSED", "chrome"),
    @(346, "Create Nationality", "PASSED", "chrome"),
    @(347, "Fee Functionality", "FAILED", "chrome"),
    @(348, "Fee Functionality", "FAILED", "chrome"),
    @(349, "Fee Functionality", "FAILED", "chrome"),
    @(350, "Fee Functionality", "FAILED", "chrome"),
    @(351, "Fee Functionality", "FAILED", "chrome"),
    @(352, "Fee Functionality", "FAILED", "chrome"),
    @(353, "Fee Functionality", "FAILED", "chrome"),
    @(354, "Fee Functionality", "FAILED", "chrome"),
    @(355, "Fee Functionality", "FAILED", "chrome"),
    @(356, "Fee Functionality", "PASSED", "chrome"),
    @(357, "Create Country", "PASSED", "chrome"),
    @(358, "Create Nationality", "PASSED", "chrome"),
    @(359, "Fee Functionality", "PASSED", "chrome"),
    @(360, "Create Country", "PASSED", "chrome"),
    @(361, "Create Nationality", "PASSED", "chrome"),
    @(362, "Fee Functionality", "PASSED", "chrome"),
    @(363, "Create and Delete Cities", "FAILED", "chrome"),
    @(364, "Create and Delete Cities", "FAILED", "chrome"),
    @(365, "Create and Delete Cities", "PASSED", "chrome")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value() = $row[1]
    $ws.Cells.Item($r, 2).Value() = $row[2]
    $ws.Cells.Item($r, 3).Value() = $row[3]
}

Write-Host "Added rows 297-365"
